# Update auto scs (lamda_1 / lamda_2), time in ms, and auto capacity
# for the Poisson compose dictionary table on Sheet1.
#
# Column A: client index (0-based), unchanged sequence but extended to 57 rows
# Column B: lamda_1 -> now constant 33.94444444444444 (was 8.159722222222223)
# Column C: lamda_2 -> now constant 1.95 (was 1.875)
# Column D: dic_nbre_clients_poisson_2_keys -> new key sequence (5 new rows appended)
# Column E: dic_nbre_clients_prob_poisson_2_values -> new probability values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lamda1 = 33.94444444444444
$lamda2 = 1.95

$dKeys = @(0,1,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,52,53,54,58,62,68,73)
$eVals = @(0.133,0.001,0.004,0.009000000000000001,0.027,0.027,0.038,0.062,0.043,0.034,0.032,0.022,0.031,0.041,0.032,0.034,0.028,0.032,0.026,0.029,0.03,0.026,0.026,0.02,0.029,0.027,0.017,0.011,0.018,0.015,0.009000000000000001,0.01,0.004,0.006,0.009000000000000001,0.009000000000000001,0.004,0.007,0.005,0.003,0.004,0.005,0.003,0.001,0.001,0.001,0.002,0.002,0.001,0.001,0.001,0.001,0.002,0.001,0.001,0.001,0.001)

$rowCount = $dKeys.Length

$data = New-Object 'object[,]' $rowCount,5
for ($i = 0; $i -lt $rowCount; $i++) {
    $data[$i,0] = $i
    $data[$i,1] = $lamda1
    $data[$i,2] = $lamda2
    $data[$i,3] = $dKeys[$i]
    $data[$i,4] = $eVals[$i]
}

$lastRow = 1 + $rowCount
$ws.Range("A2:E$lastRow").Value = $data

# New rows appended past the original range (rows 54-58) need the same
# cell style (bold, boxed border, centered) applied to column A as the
# rest of the table, since brand new cells don't inherit formatting.
$origLastRow = 53
if ($lastRow -gt $origLastRow) {
    $ws.Range("A2").Copy() | Out-Null
    $ws.Range("A" + ($origLastRow + 1) + ":A" + $lastRow).PasteSpecial(-4122) | Out-Null
}

